$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''315.79'
$ws.Range("E2").Value = '''3.29%'
$ws.Range("D3").Value = '''39.46'
$ws.Range("E3").Value = '''3.07%'
$ws.Range("D4").Value = '''5.109'
$ws.Range("E4").Value = '''0.45%'
$ws.Range("D5").Value = '''0.08174'
$ws.Range("E5").Value = '''1.37%'
$ws.Range("D6").Value = '''2.038'
$ws.Range("E6").Value = '''4.68%'
$ws.Range("D7").Value = '''8.257'
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").Value = '''4.276'
$ws.Range("E8").Value = '''2.24%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '''0.9329'
$ws.Range("E9").Value = '''0.43%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.1407'
$ws.Range("E10").Value = '''-2.33%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1993'
$ws.Range("E11").Value = '''3.49%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.09119'
$ws.Range("E12").Value = '''1.05%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03527'
$ws.Range("E13").Value = '''0.29%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09816'
$ws.Range("E14").Value = '''0.31%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001401'
$ws.Range("E15").Value = '''0.31%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.006211'
$ws.Range("E16").Value = '''2.10%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.657'
$ws.Range("E17").Value = '''-1.75%'
$ws.Range("D18").Value = '''3.290'
$ws.Range("E18").Value = '''-4.63%'
$ws.Range("D19").Value = '''0.3460'
$ws.Range("E19").Value = '''-0.06%'
$ws.Range("D20").Value = '''0.1304'
$ws.Range("E20").Value = '''-0.53%'
$ws.Range("D21").Value = '''4.892'
$ws.Range("E21").Value = '''2.17%'
$ws.Range("D23").Value = '''0.04336'
$ws.Range("E23").Value = '''-0.73%'
$ws.Range("D24").Value = '''0.001226'
$ws.Range("E24").Value = '''-0.44%'
$ws.Range("D25").Value = '''0.004779'
$ws.Range("E25").Value = '''16.02%'
$ws.Range("D26").Value = '''0.0001301'
$ws.Range("E26").Value = '''-0.06%'
$ws.Range("D27").Value = '''0.0004003'
$ws.Range("E27").Value = '''-10.00%'
$ws.Range("D39").Value = '''0.02235'
$ws.Range("E39").Value = '''7.99%'
$ws.Range("D40").Value = '''0.05257'
$ws.Range("E40").Value = '''4.39%'
$ws.Range("E41").Value = '''0.82%'
$ws.Range("D42").Value = '''0.009767'
$ws.Range("E42").Value = '''-3.41%'
$ws.Range("D43").Value = '''0.1377'
$ws.Range("E43").Value = '''2.25%'
$ws.Range("D44").Value = '''0.002151'
$ws.Range("E44").Value = '''0.42%'
$ws.Range("D45").Value = '''0.009491'
$ws.Range("E45").Value = '''6.52%'
$ws.Range("D46").Value = '''0.00006451'
$ws.Range("E46").Value = '''4.49%'
$ws.Range("E47").Value = '''-0.05%'
$ws.Range("B48").Value = 'CoinbaseStockToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D48").Value = '''0.001200'
$ws.Range("E48").Value = '''-25.04%'
$ws.Range("B49").Value = 'BOLO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D49").Value = '''0.002769'
$ws.Range("E49").Value = '''-1.44%'
$ws.Range("E50").Value = '''-0.05%'
$ws.Range("E51").Value = '''-0.05%'
